# Add Data for CycleSort
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the Cycle Sort trial data (rows 12-16, columns P-U)
$data = @(
    @(12, 4, 159, 113354),
    @(13, 1, 220, 155473),
    @(14, 2, 226, 153142),
    @(15, 1, 228, 147810),
    @(16, 1, 223, 149324)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("P$r").Value = $row[1]
    $ws.Range("Q$r").Value = $row[2]
    $ws.Range("R$r").Value = $row[3]
    $ws.Range("S$r").Value = "Unmeasureable"
    $ws.Range("T$r").Value = "Unmeasureable"
    $ws.Range("U$r").Value = "Unmeasureable"
}

# Update the selection state to match the saved workbook
$ws.Range("T21").Select()
